$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Directory" column values (D2:D5): the old mixed-case path is
# replaced with the lower-cased, current repo-relative path.
$ws.Range("D2:D5").Value = "data/multimedia/documentation/"

# Move the active selection to where the author last left the cursor.
$ws.Range("D8").Select() | Out-Null
